# Actualización automática 2025-09-08 09:55:08
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": PORCELANATO total for F.V - AREA ANDINA S.A. (row 12)
# ---------------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M12").Value = 17352.47

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL": septiembre column for F.V - AREA ANDINA S.A. (row 12)
# and the septiembre totals row (row 24)
# ---------------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F12").Value = 17352.47
$wsMensual.Range("F24").Value = 18416.84

# ---------------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL": the group breakdown shrank from 17 to 13
# groups (GRANITO, LED, PANELES PU and PANELES PVC disappeared) and every
# remaining group/total got refreshed figures.
# ---------------------------------------------------------------------------
$wsCump = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Remove the obsolete category rows (delete bottom-up so row numbers above
# the deletion point stay valid): PANELES PVC (14), PANELES PU (13),
# LED (9), GRANITO (5).
$wsCump.Rows("14").Delete()
$wsCump.Rows("13").Delete()
$wsCump.Rows("9").Delete()
$wsCump.Rows("5").Delete()

# Refresh the column F width (25 -> 24 characters stored width).
$wsCump.Columns("F").ColumnWidth = 23.1667

# Refresh every remaining category row with the new figures.
$wsCump.Range("C2").Value = 440.717086537713
$wsCump.Range("D2").Value = 0
$wsCump.Range("E2").Value = 440.717086537713
$wsCump.Range("F2").Value = 0

$wsCump.Range("C3").Value = 7475.4083879616
$wsCump.Range("D3").Value = 0
$wsCump.Range("E3").Value = 7475.4083879616
$wsCump.Range("F3").Value = 0

$wsCump.Range("C4").Value = 485.098641648355
$wsCump.Range("D4").Value = 0
$wsCump.Range("E4").Value = 485.098641648355
$wsCump.Range("F4").Value = 0

$wsCump.Range("B5").Value = "GRIFERIAS"
$wsCump.Range("C5").Value = 150
$wsCump.Range("D5").Value = 0
$wsCump.Range("E5").Value = 150
$wsCump.Range("F5").Value = 0

$wsCump.Range("B6").Value = "INODOROS"
$wsCump.Range("C6").Value = 907.166108615601
$wsCump.Range("D6").Value = 306.24
$wsCump.Range("E6").Value = 600.926108615601
$wsCump.Range("F6").Value = 0.3375787489099915

$wsCump.Range("B7").Value = "LAVABOS"
$wsCump.Range("C7").Value = 665.033262215681
$wsCump.Range("D7").Value = 0
$wsCump.Range("E7").Value = 665.033262215681
$wsCump.Range("F7").Value = 0

$wsCump.Range("B8").Value = "NO RESURTIBLES"
$wsCump.Range("C8").Value = 666.586827568148
$wsCump.Range("D8").Value = 0
$wsCump.Range("E8").Value = 666.586827568148
$wsCump.Range("F8").Value = 0

$wsCump.Range("B9").Value = "OTROS"
$wsCump.Range("C9").Value = 0
$wsCump.Range("D9").Value = 0
$wsCump.Range("E9").Value = 0
$wsCump.Range("F9").Value = 0

$wsCump.Range("B10").Value = "PANELES DECORATIVOS"
$wsCump.Range("C10").Value = 388.107983534392
$wsCump.Range("D10").Value = 0
$wsCump.Range("E10").Value = 388.107983534392
$wsCump.Range("F10").Value = 0

$wsCump.Range("B11").Value = "PIEDRA SINTERIZADA"
$wsCump.Range("C11").Value = 2922.22458185274
$wsCump.Range("D11").Value = 0
$wsCump.Range("E11").Value = 2922.22458185274
$wsCump.Range("F11").Value = 0

$wsCump.Range("B12").Value = "PORCELANATO"
$wsCump.Range("C12").Value = 43100.0854117774
$wsCump.Range("D12").Value = 18110.6
$wsCump.Range("E12").Value = 24989.4854117774
$wsCump.Range("F12").Value = 0.42019870325016

$wsCump.Range("B13").Value = "PUERTAS DE SEGURIDAD"
$wsCump.Range("C13").Value = 148.058220160454
$wsCump.Range("D13").Value = 0
$wsCump.Range("E13").Value = 148.058220160454
$wsCump.Range("F13").Value = 0

$wsCump.Range("B14").Value = "SAL SOLUBLE"
$wsCump.Range("C14").Value = 854.979720622497
$wsCump.Range("D14").Value = 0
$wsCump.Range("E14").Value = 854.979720622497
$wsCump.Range("F14").Value = 0

# Row 15 is now the TOTAL row (previously row 19).
$wsCump.Range("C15").Value = 58203.46623249458
$wsCump.Range("D15").Value = 18416.84
$wsCump.Range("E15").Value = 39786.62623249458
$wsCump.Range("F15").Value = 0.3164217046186505

Write-Host "Applied update to VENTAS POR GRUPO, VENTA MENSUAL and CUMPLIMIENTO MENSUAL"
